# Swap the order of "System" and the email address in column G wherever
# the cell currently reads "System, dnasr281@gmail.com", turning it into
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G is the 7th column
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
